# Add a new conference talk entry at the top of the data table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row right after the header row (row 1), pushing all
# existing talk rows down by one.
$ws.Rows("2:2").Insert()

# Populate the new row with the new talk's details.
$ws.Range("A2").Value = "针对中国年轻人控烟广告有潜力的信念"
$ws.Range("B2").Value = "裴瑞，于莲，赵亮，陈静茜"
$ws.Range("C2").Value = 2020
$ws.Range("D2").Value = "November"
$ws.Range("E2").Value = "The Medicine, Humanity and Media: Health China & Health Communication."
$ws.Range("F2").Value = "Online/Beijing"
$ws.Range("H2").Value = "https://www.bilibili.com/video/BV1Uv4116737"

# Reflect the cell that was last selected when the workbook was saved.
$ws.Range("H12").Select() | Out-Null
